$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H74").Value = 4319.5
$ws.Range("J74").Value = 4366.727
$ws.Range("L74").Value = 4366.727
$ws.Range("N74").Value = -6238.727
$ws.Range("H76").Value = 3211.7646
$ws.Range("I76").Value = 3257.1428
$ws.Range("J76").Value = 3000
$ws.Range("K76").Value = 3257.1428
$ws.Range("L76").Value = 3000
$ws.Range("M76").Value = -2942.1428
$ws.Range("N76").Value = -3630
$ws.Range("H77").Value = 4319.5
$ws.Range("J77").Value = 4366.727
$ws.Range("L77").Value = 21833.635
$ws.Range("N77").Value = -31193.635
$ws.Range("H79").Value = 3211.7646
$ws.Range("I79").Value = 3257.1428
$ws.Range("J79").Value = 3000
$ws.Range("K79").Value = 3257.1428
$ws.Range("L79").Value = 3000
$ws.Range("M79").Value = -2165.1428
$ws.Range("N79").Value = -5184
$ws.Range("H105").Value = 70335.5
$ws.Range("J105").Value = 70335.5
$ws.Range("L105").Value = 70335.5
$ws.Range("N105").Value = -77323.5
$ws.Range("H107").Value = 1449.75
$ws.Range("I107").Value = 0
$ws.Range("J107").Value = 1449.75
$ws.Range("K107").Value = 0
$ws.Range("L107").Value = 1449.75
$ws.Range("N107").Value = -5289.75
$ws.Range("M107").ClearContents()
$ws.Range("H112").Value = 6131.902
$ws.Range("J112").Value = 6216.54
$ws.Range("L112").Value = 18649.62
$ws.Range("N112").Value = -20865.62
$ws.Range("H115").Value = 1719.5714
$ws.Range("I115").Value = 1719.5714
$ws.Range("K115").Value = 5158.7142
$ws.Range("M115").Value = -3591.7142
$ws.Range("H118").Value = 2331.4546
$ws.Range("I118").Value = 3154
$ws.Range("J118").Value = 1646
$ws.Range("K118").Value = 9462
$ws.Range("L118").Value = 4938
$ws.Range("M118").Value = -7805
$ws.Range("N118").Value = -8252
$ws.Range("H123").Value = 0
$ws.Range("I123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("K123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("M123").ClearContents()
$ws.Range("N123").ClearContents()
$ws.Range("H128").Value = 0
$ws.Range("J128").Value = 0
$ws.Range("L128").Value = 0
$ws.Range("N128").ClearContents()
$ws.Range("H132").Value = 4569.25
$ws.Range("I132").Value = 4210.567
$ws.Range("K132").Value = 12631.701
$ws.Range("M132").Value = -10101.701
$ws.Range("H133").Value = 98780
$ws.Range("J133").Value = 98780
$ws.Range("L133").Value = 98780
$ws.Range("N133").Value = -108900
$ws.Range("H137").Value = 3448.2856
$ws.Range("I137").Value = 3031.7585
$ws.Range("J137").Value = 5461.5
$ws.Range("K137").Value = 9095.2755
$ws.Range("L137").Value = 16384.5
$ws.Range("N137").Value = -21484.5
$ws.Range("M137").Value = -6545.2755
$ws.Range("H138").Value = 1555.3662
$ws.Range("I138").Value = 2312.0527
$ws.Range("J138").Value = 1278.8846
$ws.Range("K138").Value = 6936.158100000001
$ws.Range("L138").Value = 3836.6538
$ws.Range("M138").Value = -1796.158100000001
$ws.Range("N138").Value = -14116.6538
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H46").Value = 1850
$ws.Range("J46").Value = 1850
$ws.Range("L46").Value = 1850
$ws.Range("N46").Value = -2488
$ws.Range("H63").Value = 5009.8716
$ws.Range("I63").Value = 3720
$ws.Range("J63").Value = 5199.5586
$ws.Range("K63").Value = 3720
$ws.Range("L63").Value = 5199.5586
$ws.Range("M63").Value = -3034
$ws.Range("N63").Value = -6571.5586
$ws.Range("H66").Value = 5009.8716
$ws.Range("I66").Value = 3720
$ws.Range("J66").Value = 5199.5586
$ws.Range("K66").Value = 18600
$ws.Range("L66").Value = 25997.793
$ws.Range("M66").Value = -15168
$ws.Range("N66").Value = -32861.79300000001
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82").Value = 14957.353
$ws.Range("I82").Value = 5035.5
$ws.Range("J82").Value = 23776.777
$ws.Range("K82").Value = 5035.5
$ws.Range("L82").Value = 23776.777
$ws.Range("M82").Value = -4652.5
$ws.Range("N82").Value = -24542.777
$ws.Range("H85").Value = 14957.353
$ws.Range("I85").Value = 5035.5
$ws.Range("J85").Value = 23776.777
$ws.Range("K85").Value = 5035.5
$ws.Range("L85").Value = 23776.777
$ws.Range("M85").Value = -3709.5
$ws.Range("N85").Value = -26428.777
$ws.Range("H86").Value = 52634140
$ws.Range("I86").Value = 66669308
$ws.Range("K86").Value = 66669308
$ws.Range("M86").Value = -66668185
$ws.Range("H89").Value = 52634140
$ws.Range("I89").Value = 66669308
$ws.Range("K89").Value = 333346540
$ws.Range("M89").Value = -333340924
$ws.Range("H94").Value = 956.63635
$ws.Range("I94").Value = 674
$ws.Range("K94").Value = 674
$ws.Range("M94").Value = -223
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H39").Value = 1569.9286
$ws.Range("J39").Value = 1569.9286
$ws.Range("L39").Value = 4709.7858
$ws.Range("N39").Value = -5297.7858
$ws.Range("H113").Value = 717.8372000000001
$ws.Range("I113").Value = 447.0645
$ws.Range("J113").Value = 1417.3334
$ws.Range("K113").Value = 1341.1935
$ws.Range("L113").Value = 4252.0002
$ws.Range("M113").Value = 828.8064999999999
$ws.Range("N113").Value = -8592.0002
$ws.Range("H134").Value = 5999.353
$ws.Range("J134").Value = 6636.273
$ws.Range("L134").Value = 19908.819
$ws.Range("N134").Value = -30048.819
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H18").Value = 8500
$ws.Range("J18").Value = 8500
$ws.Range("L18").Value = 8500
$ws.Range("N18").Value = -9086
$ws.Range("H70").Value = 5402.3555
$ws.Range("I70").Value = 5407.0713
$ws.Range("J70").Value = 5394.5884
$ws.Range("K70").Value = 5407.0713
$ws.Range("L70").Value = 5394.5884
$ws.Range("M70").Value = -5137.0713
$ws.Range("N70").Value = -5934.5884
$ws.Range("H73").Value = 5402.3555
$ws.Range("I73").Value = 5407.0713
$ws.Range("J73").Value = 5394.5884
$ws.Range("K73").Value = 5407.0713
$ws.Range("L73").Value = 5394.5884
$ws.Range("M73").Value = -4471.0713
$ws.Range("N73").Value = -7266.5884
$ws.Range("H97").Value = 3035.1667
$ws.Range("I97").Value = 2666.6667
$ws.Range("J97").Value = 3403.6667
$ws.Range("K97").Value = 2666.6667
$ws.Range("L97").Value = 3403.6667
$ws.Range("M97").Value = -2170.6667
$ws.Range("N97").Value = -4395.6667
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2659
$ws.Range("I40").Value = 1584.2858
$ws.Range("J40").Value = 5166.6665
$ws.Range("K40").Value = 1584.2858
$ws.Range("L40").Value = 5166.6665
$ws.Range("M40").Value = -1448.2858
$ws.Range("N40").Value = -5438.6665
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 3040.4
$ws.Range("I113").Value = 4534
$ws.Range("J113").Value = 800
$ws.Range("K113").Value = 13602
$ws.Range("L113").Value = 2400
$ws.Range("M113").Value = -11432
$ws.Range("N113").Value = -6740
$ws.Range("H122").Value = 3900.2
$ws.Range("I122").Value = 2071.8572
$ws.Range("J122").Value = 5500
$ws.Range("K122").Value = 6215.571599999999
$ws.Range("L122").Value = 16500
$ws.Range("M122").Value = -3765.571599999999
$ws.Range("N122").Value = -21400
